# Adds two new attribute rows ("foto" and "tipo_usuario") to the end of
# the tb_usuario table (the first table in the document, which currently
# ends with the "senha : varchar (255)" / "Senha do usuário" row).

$d = $word.ActiveDocument

# The tb_usuario table is the first table in the document.
$table = $d.Tables.Item(1)

# Row: foto: varchar | Foto do usuário
$rowFoto = $table.Rows.Add()
$rowFoto.Cells.Item(1).Range.Text = "foto: varchar"
$rowFoto.Cells.Item(2).Range.Text = "Foto do usuário"

# Row: tipo_usuario | Tipo do usuário: Administrador ou Comum
$rowTipo = $table.Rows.Add()
$rowTipo.Cells.Item(1).Range.Text = "tipo_usuario"
$rowTipo.Cells.Item(2).Range.Text = "Tipo do usuário: Administrador ou Comum "
